$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 1.237070695960085
    "C2" = 0.3105736028194883
    "E2" = 0.7658344084233022
    "F2" = 2.080998279908741
    "G2" = 0.00239708178956788
    "J2" = 0.02995519906780153
    "O2" = 1.537663736486479
    "B3" = 1.090371711796365
    "C3" = 0.2718714026372027
    "E3" = 0.7354788959205223
    "F3" = 2.054911353883298
    "G3" = 0.002400078095524242
    "J3" = 0.03062390510857327
    "O3" = 1.571178246654483
    "B4" = 1.00006562515108
    "C4" = 0.248008428769765
    "E4" = 0.7171059318950626
    "F4" = 2.040415005971994
    "G4" = 0.002402012034680739
    "J4" = 0.03107427839695731
    "O4" = 1.593784093427303
    "B5" = 0.9632086851536883
    "C5" = 0.238259623645348
    "E5" = 0.7096856415208634
    "F5" = 2.034888857923065
    "G5" = 0.002402823890812928
    "J5" = 0.03126776277949794
    "O5" = 1.603503668571875
    "B6" = 0.9570852651589803
    "C6" = 0.2366393835161489
    "E6" = 0.7084575479520652
    "F6" = 2.033994233728876
    "G6" = 0.002402960136266938
    "J6" = 0.03130049040405147
    "O6" = 1.605148178621008
    "B7" = 0.9995687856412019
    "C7" = 0.2478770511481514
    "E7" = 0.7170055884958231
    "F7" = 2.040338936501584
    "G7" = 0.002402022887558193
    "J7" = 0.03107684756012752
    "O7" = 1.593913123323361
    "B8" = 1.186538213996869
    "C8" = 0.2972501227826285
    "E8" = 0.7553127787438285
    "F8" = 2.071687017642788
    "G8" = 0.002398095409832397
    "J8" = 0.03017748108841545
    "O8" = 1.548797230153014
    "B9" = 1.551275112127598
    "C9" = 0.3932578202549735
    "E9" = 0.8325398722141557
    "F9" = 2.145293681163267
    "G9" = 0.002391137663308026
    "J9" = 0.02873176708551028
    "O9" = 1.476521855829887
    "B10" = 1.818022220607588
    "C10" = 0.4632770363875807
    "E10" = 0.8905708947053341
    "F10" = 2.20686667647189
    "G10" = 0.002386474625608636
    "J10" = 0.02786652987599325
    "O10" = 1.433441685565811
    "B11" = 1.939095710745164
    "C11" = 0.4950140990478076
    "E11" = 0.9172535509981543
    "F11" = 2.236526809804872
    "G11" = 0.002384449729469188
    "J11" = 0.02751634026443028
    "O11" = 1.416050724877138
    "B12" = 1.984902663924061
    "C12" = 0.5070150762532535
    "E12" = 0.927398470928253
    "F12" = 2.247997225705575
    "G12" = 0.00238369673117298
    "J12" = 0.02739003049035915
    "O12" = 1.409785104465612
    "B13" = 1.975039169127683
    "C13" = 0.5044312239982105
    "E13" = 0.9252117683645054
    "F13" = 2.245516223805964
    "G13" = 0.002383858290735774
    "J13" = 0.02741695248293041
    "O13" = 1.411120245008732
    "B14" = 1.942865108176022
    "C14" = 0.4960017740731928
    "E14" = 0.9180873628056219
    "F14" = 2.237465694172982
    "G14" = 0.002384387503921934
    "J14" = 0.02750582210653008
    "O14" = 1.415528820717967
    "B15" = 1.923152176246106
    "C15" = 0.4908362397363248
    "E15" = 0.9137287671274095
    "F15" = 2.232565650369651
    "G15" = 0.002384713455667828
    "J15" = 0.02756107938237307
    "O15" = 1.418270948232902
    "B16" = 1.810104159121295
    "C16" = 0.4612005768909739
    "E16" = 0.8888328301470807
    "F16" = 2.2049616492383
    "G16" = 0.002386608890234933
    "J16" = 0.02789029399976606
    "O16" = 1.434622853526918
    "B17" = 1.740682015589812
    "C17" = 0.4429901521909301
    "E17" = 0.8736326804129817
    "F17" = 2.188451197956624
    "G17" = 0.002387796305318123
    "J17" = 0.02810341562321383
    "O17" = 1.445221241017563
    "B18" = 1.700726832038129
    "C18" = 0.4325052066733406
    "E18" = 0.8649166964575414
    "F18" = 2.17911006621199
    "G18" = 0.002388488347105329
    "J18" = 0.02823008019242401
    "O18" = 1.451524672147627
    "B19" = 1.687194390929051
    "C19" = 0.4289533499531899
    "E19" = 0.8619702089242338
    "F19" = 2.175973938348818
    "G19" = 0.002388724220988889
    "J19" = 0.02827366610435078
    "O19" = 1.45369446564851
    "B20" = 1.748074764720002
    "C20" = 0.4449298043209069
    "E20" = 0.8752479953232495
    "F20" = 2.190192686044242
    "G20" = 0.002387668964553576
    "J20" = 0.02808030551211083
    "O20" = 1.44407152773725
    "B21" = 1.952316536638477
    "C21" = 0.498478178431526
    "E21" = 0.9201788672565954
    "F21" = 2.239823838110453
    "G21" = 0.002384231687594993
    "J21" = 0.02747954752453907
    "O21" = 1.414225210223677
    "B22" = 2.085560767529671
    "C22" = 0.5333748375993537
    "E22" = 0.9497814422632018
    "F22" = 2.273653067973299
    "G22" = 0.002382065552651948
    "J22" = 0.02712366190144522
    "O22" = 1.396585125904636
    "B23" = 2.014468374837179
    "C23" = 0.5147592127319172
    "E23" = 0.9339602667164826
    "F23" = 2.255469881048668
    "G23" = 0.002383214332360729
    "J23" = 0.02731022377661319
    "O23" = 1.405828327454799
    "B24" = 1.744732638094206
    "C24" = 0.4440529359251286
    "E24" = 0.8745176404789703
    "F24" = 2.189404889189689
    "G24" = 0.002387726506040584
    "J24" = 0.0280907407080182
    "O24" = 1.444590658294999
    "B25" = 1.452814992071467
    "C25" = 0.3673745380013997
    "E25" = 0.8114215560537588
    "F25" = 2.124072149015461
    "G25" = 0.002392940758420754
    "J25" = 0.02908851292576387
    "O25" = 1.494325138371281
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}

Write-Output "Updated $($data.Count) cells"